# Update diversity metric / richness summary statistics (Schloss, Young, Jackson,
# Charles River, Taconic, Envigo columns) per updated Kruskal-Wallis analysis
# (plot medians instead of mean; commit: "Figure edits based on Nick's feedback...")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.263265
$ws.Range("G2").Value = 1.060187
$ws.Range("H2").Value = 0.88605
$ws.Range("I2").Value = 2.3578725
$ws.Range("J2").Value = 1.600417
$ws.Range("K2").Value = 1.63879
$ws.Range("F3").Value = 2.0133175
$ws.Range("G3").Value = 2.230369
$ws.Range("H3").Value = 1.350626
$ws.Range("I3").Value = 3.427412
$ws.Range("J3").Value = 2.751541
$ws.Range("K3").Value = 3.063784
$ws.Range("F4").Value = 2.985855
$ws.Range("G4").Value = 2.648099
$ws.Range("H4").Value = 2.390782
$ws.Range("I4").Value = 3.5449285
$ws.Range("J4").Value = 2.67196
$ws.Range("K4").Value = 2.546703
$ws.Range("F5").Value = 1.674286
$ws.Range("G5").Value = 1.211476
$ws.Range("H5").Value = 1.527338
$ws.Range("I5").Value = 2.937397
$ws.Range("J5").Value = 2.352139
$ws.Range("K5").Value = 1.964147
$ws.Range("F6").Value = 0.393113
$ws.Range("G6").Value = 0.6393135000000001
$ws.Range("H6").Value = 0.8242039999999999
$ws.Range("I6").Value = 1.47545
$ws.Range("J6").Value = 1.287899
$ws.Range("K6").Value = 1.1060325
$ws.Range("F7").Value = 2.7673705
$ws.Range("G7").Value = 2.577549
$ws.Range("H7").Value = 2.239505
$ws.Range("I7").Value = 3.541092
$ws.Range("J7").Value = 2.5939365
$ws.Range("K7").Value = 2.713827
$ws.Range("F8").Value = 2.0671315
$ws.Range("G8").Value = 1.849646
$ws.Range("H8").Value = 1.340242
$ws.Range("I8").Value = 3.256496
$ws.Range("J8").Value = 2.109229
$ws.Range("K8").Value = 2.411948
$ws.Range("F9").Value = 2.709545
$ws.Range("G9").Value = 3.3462055
$ws.Range("H9").Value = 2.676745
$ws.Range("I9").Value = 3.6113935
$ws.Range("J9").Value = 2.5004215
$ws.Range("K9").Value = 2.942861
$ws.Range("F10").Value = 2.010388
$ws.Range("G10").Value = 1.5602365
$ws.Range("K10").Value = 1.853689
$ws.Range("F11").Value = 2.9433485
$ws.Range("G11").Value = 2.3287495
$ws.Range("H11").Value = 2.0863455
$ws.Range("J11").Value = 2.869485
$ws.Range("K11").Value = 3.0589595
$ws.Range("F12").Value = 2.649149
$ws.Range("G12").Value = 2.68211
$ws.Range("H12").Value = 2.01088
$ws.Range("I12").Value = 3.301791
$ws.Range("J12").Value = 2.5972915
$ws.Range("K12").Value = 2.776341
$ws.Range("F13").Value = 86.5
$ws.Range("G13").Value = 104
$ws.Range("H13").Value = 59
$ws.Range("I13").Value = 171
$ws.Range("J13").Value = 89.5
$ws.Range("K13").Value = 110.5
$ws.Range("F14").Value = 115.5
$ws.Range("G14").Value = 130
$ws.Range("H14").Value = 91.5
$ws.Range("I14").Value = 191.5
$ws.Range("J14").Value = 97
$ws.Range("K14").Value = 126
$ws.Range("F15").Value = 104.5
$ws.Range("G15").Value = 127
$ws.Range("H15").Value = 80
$ws.Range("I15").Value = 197
$ws.Range("J15").Value = 93.5
$ws.Range("K15").Value = 110
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 26
$ws.Range("H16").Value = 15
$ws.Range("I16").Value = 80.5
$ws.Range("J16").Value = 26
$ws.Range("K16").Value = 34.5
$ws.Range("F17").Value = 61
$ws.Range("G17").Value = 55.5
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 146.5
$ws.Range("J17").Value = 61
$ws.Range("K17").Value = 73
$ws.Range("F18").Value = 111
$ws.Range("G18").Value = 126
$ws.Range("H18").Value = 87
$ws.Range("I18").Value = 176
$ws.Range("J18").Value = 86
$ws.Range("K18").Value = 119
$ws.Range("F19").Value = 83
$ws.Range("G19").Value = 86
$ws.Range("H19").Value = 53
$ws.Range("I19").Value = 155
$ws.Range("J19").Value = 85
$ws.Range("K19").Value = 96
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 18.5
$ws.Range("H20").Value = 14.5
$ws.Range("I20").Value = 34
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = 21
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 131.5
$ws.Range("H21").Value = 134
$ws.Range("I21").Value = 247
$ws.Range("J21").Value = 122.5
$ws.Range("K21").Value = 142
$ws.Range("F22").Value = 111.5
$ws.Range("G22").Value = 109
$ws.Range("H22").Value = 85
$ws.Range("J22").Value = 105
$ws.Range("K22").Value = 124
$ws.Range("F23").Value = 42
$ws.Range("G23").Value = 45.5
$ws.Range("K23").Value = 51
